$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell used as a format donor (plain/default style, no special formatting)
$formatDonor = "B2"

$ws.Range('D2').Value = '41.829.00'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '2.479.38'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range($formatDonor).Copy()
$ws.Range('D4').PasteSpecial(-4122)
$ws.Range('D5').Value = '''319.16'
$ws.Range($formatDonor).Copy()
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').Value = '''93.52'
$ws.Range($formatDonor).Copy()
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.518'
$ws.Range($formatDonor).Copy()
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '''0.0879'
$ws.Range($formatDonor).Copy()
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('E10').Value = '  +11.25%  '
$ws.Range('D11').Value = '''33.36'
$ws.Range($formatDonor).Copy()
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('E11').Value = '  +2.62%  '
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').Value = '2.861.29'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').Value = '''6.97'
$ws.Range($formatDonor).Copy()
$ws.Range('D14').PasteSpecial(-4122)
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').Value = '''15.72'
$ws.Range($formatDonor).Copy()
$ws.Range('D15').PasteSpecial(-4122)
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').Value = '2.460.75'
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').Value = '''0.800'
$ws.Range($formatDonor).Copy()
$ws.Range('D17').PasteSpecial(-4122)
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('D18').Value = '41.781.62'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').Value = '0.0₃0955'
$ws.Range('E19').Value = '  +1.54%  '
$ws.Range('D20').Value = '''6.48'
$ws.Range($formatDonor).Copy()
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').Value = '''71.28'
$ws.Range($formatDonor).Copy()
$ws.Range('D21').PasteSpecial(-4122)
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('D22').Value = '''11.40'
$ws.Range($formatDonor).Copy()
$ws.Range('D22').PasteSpecial(-4122)
$ws.Range('E22').Value = '  +2.85%  '
$ws.Range('D23').Value = '''242.18'
$ws.Range($formatDonor).Copy()
$ws.Range('D23').PasteSpecial(-4122)
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('D25').Value = '''1.97'
$ws.Range($formatDonor).Copy()
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = '''25.38'
$ws.Range($formatDonor).Copy()
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('E27').Value = '  +3.17%  '
$ws.Range('D28').Value = '''2.26'
$ws.Range($formatDonor).Copy()
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').Value = '''9.81'
$ws.Range($formatDonor).Copy()
$ws.Range('D29').PasteSpecial(-4122)
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('D30').Value = '''37.25'
$ws.Range($formatDonor).Copy()
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('E30').Value = '  +5.14%  '
$ws.Range('D31').Value = '''157.65'
$ws.Range($formatDonor).Copy()
$ws.Range('D31').PasteSpecial(-4122)
$ws.Range('E31').Value = '  +1.22%  '
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '''0.0769'
$ws.Range($formatDonor).Copy()
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '''17.58'
$ws.Range($formatDonor).Copy()
$ws.Range('D36').PasteSpecial(-4122)
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('D37').Value = '''1.89'
$ws.Range($formatDonor).Copy()
$ws.Range('D37').PasteSpecial(-4122)
$ws.Range('E37').Value = '  +5.89%  '
$ws.Range('E38').Value = '  +2.42%  '
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('E40').Value = '  +1.31%  '
$ws.Range('E41').Value = '  +7.37%  '
$ws.Range('D42').Value = '''4.03'
$ws.Range($formatDonor).Copy()
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('E42').Value = '  +0.98%  '
$ws.Range('D43').Value = '2.007.54'
$ws.Range('E43').Value = '  +3.36%  '
$ws.Range('D44').Value = '''19.16'
$ws.Range($formatDonor).Copy()
$ws.Range('D44').PasteSpecial(-4122)
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D47').Value = '''9.48'
$ws.Range($formatDonor).Copy()
$ws.Range('D47').PasteSpecial(-4122)
$ws.Range('E47').Value = '  +4.67%  '
$ws.Range('D48').Value = '2.718.41'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = '''98.55'
$ws.Range($formatDonor).Copy()
$ws.Range('D49').PasteSpecial(-4122)
$ws.Range('E49').Value = '  +1.47%  '
$ws.Range('D50').Value = '''76.62'
$ws.Range($formatDonor).Copy()
$ws.Range('D50').PasteSpecial(-4122)
$ws.Range('E50').Value = '  +7.05%  '
$ws.Range('D51').Value = '''67.76'
$ws.Range($formatDonor).Copy()
$ws.Range('D51').PasteSpecial(-4122)
$ws.Range('E51').Value = '  +1.18%  '

$excel.CutCopyMode = 0
